$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("per image")

# Insert 3 new rows before row 20 (CK19 - A - 22/23/24)
$ws.Rows.Item(20).Resize(3).Insert(-4121)  # xlShiftDown = -4121

$ws.Range("A20").Value = "CK19 - A - 22"
$ws.Range("B20").Value = 0.666
$ws.Range("C20").Value = "ConA"

$ws.Range("A21").Value = "CK19 - A - 23"
$ws.Range("B21").Value = 0.264
$ws.Range("C21").Value = "OVA"

$ws.Range("A22").Value = "CK19 - A - 24"
$ws.Range("B22").Value = 0.193
$ws.Range("C22").Value = "Unstimulated"

# Insert 3 new rows before row 41 (CK19 - B - 22/23/24), which corresponds
# to old row 38 ("CK19 - B - 19") shifted down by 3 already.
$ws.Rows.Item(41).Resize(3).Insert(-4121)  # xlShiftDown = -4121

$ws.Range("A41").Value = "CK19 - B - 22"
$ws.Range("B41").Value = 0.465
$ws.Range("C41").Value = "ConA"

$ws.Range("A42").Value = "CK19 - B - 23"
$ws.Range("B42").Value = 0.297
$ws.Range("C42").Value = "OVA"

$ws.Range("A43").Value = "CK19 - B - 24"
$ws.Range("B43").Value = 0.209
$ws.Range("C43").Value = "Unstimulated"

$ws.Range("C41").Select()

$wb.Save()
